# Sprint 1 Backlog.xlsx - "Added task to create JUnit Tests"
#
# Adds a new backlog row (item 19 / task #16) describing the JUnit-tests
# task, gives it an estimate of 4, and matches the vertical+horizontal
# centering that the author applied to the "Est" (column D) cells for the
# surrounding rows (14-20). Also moves the saved selection to E17, matching
# where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center (horizontally + vertically) the "Est" column cells for rows 14-20,
# matching the style already used for the "Item #"/"Description" columns in
# that block (bordered, centered) but adding vertical centering too.
$estRange = $ws.Range("D14:D20")
$estRange.HorizontalAlignment = -4108   # xlCenter
$estRange.VerticalAlignment = -4108     # xlCenter

# Fill in the new backlog item: "Create Junit tests for equivalent classes",
# estimated at 4.
$ws.Range("C19").Value = "Create Junit tests for equivalent classes"
$ws.Range("D19").Value = 4

# Match the author's final cursor position.
$ws.Range("E17").Select() | Out-Null
